$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1 (was the source filename, now a generic title)
$ws.Range("A1").Value = "Microstate List"

# The four microstates that turned out to be redundant resonance
# structures / duplicate geometric isomers (micro001, micro012, micro025,
# micro035) are dropped from the table. Rather than deleting those rows
# outright (which in Excel would also drag each row's own formatting up
# with it), copy the surviving rows' *values* upward into place so that
# each row keeps its own pre-existing alternating shading - exactly
# matching how the published workbook's row bands stayed fixed to the row
# number while only the cell contents shifted.
$mapping = @{
  3=4;  4=5;  5=6;  6=7;  7=8;  8=9;  9=10;  10=11; 11=12; 12=13;
  13=15; 14=16; 15=17; 16=18; 17=19; 18=20; 19=21; 20=22; 21=23; 22=24; 23=25; 24=26;
  25=28; 26=29; 27=30; 28=31; 29=32; 30=33; 31=34; 32=35; 33=36; 34=38
}

foreach ($fr in 3..34) {
    $orr = $mapping[$fr]
    $ws.Range("B" + $fr).Value = $ws.Range("B" + $orr).Value()
    $ws.Range("C" + $fr).Value = $ws.Range("C" + $orr).Value()
}

# The last four rows (which used to hold micro033-036 before the shift)
# are no longer needed now that everything has moved up four slots.
$ws.Range("A35:A38").EntireRow.Delete()

# The picture objects are anchored independently of the row data and were
# not reshuffled to track the content shift; only the trailing four
# pictures were removed, matching the rows dropped from the end of the
# sheet.
$ws.Shapes.Item("Picture 36").Delete()
$ws.Shapes.Item("Picture 35").Delete()
$ws.Shapes.Item("Picture 34").Delete()
$ws.Shapes.Item("Picture 33").Delete()
